# Apply edits to the RDF sample workbook as described by the commit diff.
$wb = $excel.ActiveWorkbook

$wsGood = $wb.Worksheets.Item("good")
$wsBad1 = $wb.Worksheets.Item("bad")
$wsBad2 = $wb.Worksheets.Item("bad (2)")

# --- Sheet "good" updates ---
$wsGood.Range("A1").Value = "extra_header"
$wsGood.Range("C6").Value = "NO ELEMENT PROVIDED. THIS SHOULD BE IGNORED."
$wsGood.Range("A2").Value = "THIS SHOULD BE IGNORED."
$wsGood.Range("A3").Value = "THIS SHOULD BE IGNORED."
$wsGood.Range("C8").Value = "id_1234"

# Column width for column A
$wsGood.Columns.Item(1).ColumnWidth = 23.666666666666668

# Styling: header cell C1 gets bold font + left alignment; the rest of column C
# (cells with content) gets left alignment only.
$wsGood.Range("C1").HorizontalAlignment = -4131
$wsGood.Range("C2").HorizontalAlignment = -4131
$wsGood.Range("C3").HorizontalAlignment = -4131
$wsGood.Range("C4").HorizontalAlignment = -4131
$wsGood.Range("C6").HorizontalAlignment = -4131
$wsGood.Range("C7").HorizontalAlignment = -4131
$wsGood.Range("C8").HorizontalAlignment = -4131

# --- Sheet "bad" updates ---
$wsBad1.Range("A1").Value = "ignored"
$wsBad1.Range("A2").Value = "bad"
$wsBad1.Range("A3").Value = "bad"
$wsBad1.Columns.Item(1).ColumnWidth = 7.0
$wsBad1.Range("A4").Select()

# --- Sheet "bad (2)" updates ---
$wsBad2.Range("A1").Value = "ignored"
$wsBad2.Range("A2").Value = "bad"
$wsBad2.Range("A3").Value = "bad"
$wsBad2.Columns.Item(1).ColumnWidth = 7.0
$wsBad2.Range("C11:C12").Select()

# Re-activate "good" as the selected tab, with its own selection restored.
$wsGood.Range("A4").Select()
$wsGood.Activate()
